# daily auto push: 2026-02-22 09:43 UTC
#
# The source data table (A:date, B:weekday, C:hour, D:ranking) gets one new
# sample appended for 2026/02/22 09:43 UTC. In the sheet that lands as a new
# row inserted right after the existing last "2026/02/22" row (row 859),
# which pushes every following row down by one and grows the used range
# from A1:D901 to A1:D902.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at 860; rows 860-901 (and the sheet dimension) shift
# down to 861-902 automatically.
$ws.Rows(860).Insert()

# The date column holds plain text like "2026/12/29", not real dates, so
# force the new cell to Text first - otherwise Excel would happily "helpfully"
# reinterpret the "YYYY/MM/DD" string as a date serial number.
$ws.Range("A860").NumberFormat = "@"
$ws.Range("A860").Value = "2026/02/22"
$ws.Range("B860").Value = "日"
$ws.Range("C860").Value = 16
$ws.Range("D860").Value = 201

# Pull the plain (unstyled) formatting from the neighboring row back onto
# the new cell so it doesn't end up carrying a stray "Text" style index -
# every other data cell in the sheet is unstyled too.
$ws.Range("A859").Copy()
$ws.Range("A860").PasteSpecial(-4122)
